$wb = $excel.ActiveWorkbook

# Overview sheet: update Latest HO Xliff Generate Date for the second file
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-27 02:47:45"

# zh-cn sheet: update Correspond Handoff Datetime / Correspond Handback DateTime for row 3
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-27 02:47:40"
$wsZhCn.Range("K3").Value = "2016-08-27 02:47:57"

# de-de sheet: update Correspond Handoff Datetime / Correspond Handback DateTime for rows 2 and 3
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-27 02:47:45"
$wsDeDe.Range("H3").Value = "2016-08-27 02:47:45"
$wsDeDe.Range("K3").Value = "2016-08-27 02:48:09"
